# Weekly data refresh: insert two new daily-price rows at the top of the
# data block (row 38/39, right after the header + rows 1-37 that are left
# untouched) and push the previously-existing rows 38:69 down to 40:71.
#
# Column layout (row 1 header):
# A Mercado ID | B Mercado | C Región | D Fecha | E Codreg | F Categoría ID
# G Categoría | H Variedad | I Calidad | J Volumen | K Precio mínimo
# L Precio máximo | M Precio promedio ponderado | N Unidad de comercialización
# O Origen | P Precio $/Kg | Q Kg o Unidades | R Clasificación

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 38:69 down to 40:71, inserting two blank rows at 38:39.
$ws.Rows("38:39").Insert()

# New row 38
$ws.Cells.Item(38, 1).Value = 7
$ws.Cells.Item(38, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(38, 3).Value = "Ñuble"
$ws.Cells.Item(38, 4).Value = 45240
$ws.Cells.Item(38, 5).Value = 16
$ws.Cells.Item(38, 6).Value = 300000000
$ws.Cells.Item(38, 7).Value = "Espárragos"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 500
$ws.Cells.Item(38, 11).Value = 1200
$ws.Cells.Item(38, 12).Value = 1200
$ws.Cells.Item(38, 13).Value = 1200
$ws.Cells.Item(38, 14).Value = "$/kilo"
$ws.Cells.Item(38, 15).Value = "Región de Ñuble"
$ws.Cells.Item(38, 16).Value = 1200
$ws.Cells.Item(38, 17).Value = 1
$ws.Cells.Item(38, 18).Value = "Hortaliza"

# New row 39
$ws.Cells.Item(39, 1).Value = 7
$ws.Cells.Item(39, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(39, 3).Value = "Ñuble"
$ws.Cells.Item(39, 4).Value = 45240
$ws.Cells.Item(39, 5).Value = 16
$ws.Cells.Item(39, 6).Value = 300000000
$ws.Cells.Item(39, 7).Value = "Espárragos"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Segunda"
$ws.Cells.Item(39, 10).Value = 400
$ws.Cells.Item(39, 11).Value = 1000
$ws.Cells.Item(39, 12).Value = 1000
$ws.Cells.Item(39, 13).Value = 1000
$ws.Cells.Item(39, 14).Value = "$/kilo"
$ws.Cells.Item(39, 15).Value = "Región de Ñuble"
$ws.Cells.Item(39, 16).Value = 1000
$ws.Cells.Item(39, 17).Value = 1
$ws.Cells.Item(39, 18).Value = "Hortaliza"

$ws.Range("A1").Select()
